$wb = $excel.ActiveWorkbook

# Mapping of old sheet name -> new sheet name and the corresponding
# title text that lives in cell A1 of each sheet.
$renames = @(
    @{ Old = "GP1"; New = "GP01"; Title = "KPI GP01 - Global Perf 1" },
    @{ Old = "GP2"; New = "GP02"; Title = "KPI GP02 - Global Perf 2" },
    @{ Old = "BP1"; New = "BP01"; Title = "KPI BP01 - Business Process 1" },
    @{ Old = "BP2"; New = "BP02"; Title = "KPI BP02 - Business Process 2" },
    @{ Old = "BP3"; New = "BP03"; Title = "KPI BP03 - Business Process 3" },
    @{ Old = "BP4"; New = "BP04"; Title = "KPI BP04 - Business Process 4" },
    @{ Old = "BP5"; New = "BP05"; Title = "KPI BP05 - Business Process 5" },
    @{ Old = "BP6"; New = "BP06"; Title = "KPI BP06 - Business Process 6" },
    @{ Old = "BP7"; New = "BP07"; Title = "KPI BP07 - Business Process 7" }
)

foreach ($r in $renames) {
    $ws = $wb.Worksheets.Item($r.Old)
    $ws.Range("A1").Value = $r.Title
    $ws.Name = $r.New
}
